# Insert a new leading run containing a single space " " (Arial, 10pt)
# immediately before the existing first run of the document's title
# paragraph ("Identificación del problema y análisis de requerimientos"),
# so the run is split into two runs: a new one holding just " " and the
# pre-existing one (now reduced to its original content).

$d = $word.ActiveDocument

# Use Find/Replace (rather than Range.InsertBefore) so the newly created
# text does not pick up an extra revision-session rsid attribute on the
# run - it grows the leading space that will become its own run once we
# touch its formatting below.
$rng = $d.Content
$rng.Find.Execute(" Identificación del problema y análisis de requerimientos", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "  Identificación del problema y análisis de requerimientos", 2) | Out-Null

# Re-apply the (identical) Arial/10pt formatting to just the new leading
# space so Word splits it into its own run, separate from the rest of
# the title text.
$newRun = $d.Range(0, 1)
$newRun.Font.NameAscii = "Arial"
$newRun.Font.Name = "Arial"
$newRun.Font.NameBi = "Arial"
$newRun.Font.Size = 10
